$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3274.525
$ws.Range("I132").Value = 2999.3684
$ws.Range("J132").Value = 8502.5
$ws.Range("K132").Value = 8998.1052
$ws.Range("L132").Value = 25507.5
$ws.Range("M132").Value = -6468.1052
$ws.Range("N132").Value = -30567.5
$ws.Range("H135").Value = 1531.8096
$ws.Range("I135").Value = 746.75
$ws.Range("J135").Value = 4044
$ws.Range("K135").Value = 6720.75
$ws.Range("L135").Value = 36396
$ws.Range("M135").Value = -4185.75
$ws.Range("N135").Value = -41466
$ws.Range("H138").Value = 2414.5286
$ws.Range("J138").Value = 3444.1177
$ws.Range("L138").Value = 10332.3531
$ws.Range("N138").Value = -20612.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 2690
$ws.Range("I25").Value = 1900
$ws.Range("K25").Value = 1900
$ws.Range("M25").Value = -1498
$ws.Range("H74").Value = 1516.75
$ws.Range("I74").Value = 850.7778
$ws.Range("J74").Value = 2715.5
$ws.Range("K74").Value = 850.7778
$ws.Range("L74").Value = 2715.5
$ws.Range("M74").Value = 23.22220000000004
$ws.Range("N74").Value = -4463.5
$ws.Range("H77").Value = 1516.75
$ws.Range("I77").Value = 850.7778
$ws.Range("J77").Value = 2715.5
$ws.Range("K77").Value = 4253.889
$ws.Range("L77").Value = 13577.5
$ws.Range("M77").Value = 114.1109999999999
$ws.Range("N77").Value = -22313.5
$ws.Range("H129").Value = 45223
$ws.Range("J129").Value = 45223
$ws.Range("L129").Value = 45223
$ws.Range("N129").Value = -55223
$ws.Range("H132").Value = 2988.2812
$ws.Range("I132").Value = 2028.76
$ws.Range("J132").Value = 6415.143
$ws.Range("K132").Value = 6086.28
$ws.Range("L132").Value = 19245.429
$ws.Range("M132").Value = -3556.28
$ws.Range("N132").Value = -24305.429

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 50000
$ws.Range("J58").Value = 50000
$ws.Range("L58").Value = 50000
$ws.Range("N58").Value = -50588
$ws.Range("H86").Value = 3082.4285
$ws.Range("I86").Value = 3236
$ws.Range("J86").Value = 2698.5
$ws.Range("K86").Value = 3236
$ws.Range("L86").Value = 2698.5
$ws.Range("M86").Value = -2113
$ws.Range("N86").Value = -4944.5
$ws.Range("H89").Value = 3082.4285
$ws.Range("I89").Value = 3236
$ws.Range("J89").Value = 2698.5
$ws.Range("K89").Value = 16180
$ws.Range("L89").Value = 13492.5
$ws.Range("M89").Value = -10564
$ws.Range("N89").Value = -24724.5
$ws.Range("H94").Value = 1489.2
$ws.Range("I94").Value = 1148.2354
$ws.Range("J94").Value = 2213.75
$ws.Range("K94").Value = 1148.2354
$ws.Range("L94").Value = 2213.75
$ws.Range("M94").Value = -697.2354
$ws.Range("N94").Value = -3115.75
$ws.Range("H133").Value = 47833.332
$ws.Range("J133").Value = 47833.332
$ws.Range("L133").Value = 47833.332
$ws.Range("N133").Value = -57953.332
$ws.Range("H134").Value = 2041.2812
$ws.Range("I134").Value = 1761.4642
$ws.Range("K134").Value = 5284.392599999999
$ws.Range("M134").Value = -2749.392599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 832.3333
$ws.Range("I16").Value = 659.3333
$ws.Range("K16").Value = 659.3333
$ws.Range("M16").Value = -372.3333
$ws.Range("H22").Value = 770.1429000000001
$ws.Range("J22").Value = 918.2
$ws.Range("L22").Value = 918.2
$ws.Range("N22").Value = -1618.2
$ws.Range("H31").Value = 9948
$ws.Range("I31").Value = 1924
$ws.Range("J31").Value = 13768.952
$ws.Range("K31").Value = 1924
$ws.Range("L31").Value = 13768.952
$ws.Range("M31").Value = -1629
$ws.Range("N31").Value = -14358.952
$ws.Range("H34").Value = 9948
$ws.Range("I34").Value = 1924
$ws.Range("J34").Value = 13768.952
$ws.Range("K34").Value = 1924
$ws.Range("L34").Value = 13768.952
$ws.Range("M34").Value = -1722
$ws.Range("N34").Value = -14172.952
$ws.Range("H113").Value = 832.3333
$ws.Range("I113").Value = 659.3333
$ws.Range("K113").Value = 659.3333
$ws.Range("M113").Value = 1510.6667
$ws.Range("H134").Value = 5977.385
$ws.Range("I134").Value = 7077.5884
$ws.Range("J134").Value = 3899.2222
$ws.Range("K134").Value = 21232.7652
$ws.Range("L134").Value = 11697.6666
$ws.Range("M134").Value = -18697.7652
$ws.Range("N134").Value = -16767.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 807.5294
$ws.Range("J5").Value = 1597
$ws.Range("L5").Value = 4791
$ws.Range("N5").Value = -5015
$ws.Range("H9").Value = 39626.13
$ws.Range("J9").Value = 41404.59
$ws.Range("L9").Value = 124213.77
$ws.Range("N9").Value = -124661.77
$ws.Range("H12").Value = 180.84616
$ws.Range("I12").Value = 101.5
$ws.Range("J12").Value = 195.27272
$ws.Range("K12").Value = 304.5
$ws.Range("L12").Value = 585.81816
$ws.Range("M12").Value = -131.5
$ws.Range("N12").Value = -931.81816
$ws.Range("H122").Value = 3153.1025
$ws.Range("I122").Value = 348.65216
$ws.Range("J122").Value = 7184.5
$ws.Range("K122").Value = 3137.86944
$ws.Range("L122").Value = 64660.5
$ws.Range("M122").Value = -687.8694399999999
$ws.Range("N122").Value = -69560.5
$ws.Range("H131").Value = 1012.44446
$ws.Range("J131").Value = 1191.25
$ws.Range("L131").Value = 3573.75
$ws.Range("N131").Value = -13653.75
$ws.Range("H135").Value = 807.5294
$ws.Range("J135").Value = 1597
$ws.Range("L135").Value = 14373
$ws.Range("N135").Value = -19443

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 46310530
$ws.Range("I80").Value = 56556756
$ws.Range("J80").Value = 202503
$ws.Range("K80").Value = 56556756
$ws.Range("L80").Value = 202503
$ws.Range("M80").Value = -56555758
$ws.Range("N80").Value = -204499
$ws.Range("H83").Value = 46310530
$ws.Range("I83").Value = 56556756
$ws.Range("J83").Value = 202503
$ws.Range("K83").Value = 282783780
$ws.Range("L83").Value = 1012515
$ws.Range("M83").Value = -282778788
$ws.Range("N83").Value = -1022499
$ws.Range("H132").Value = 2927.6135
$ws.Range("I132").Value = 2494.5833
$ws.Range("K132").Value = 7483.749899999999
$ws.Range("M132").Value = -4953.749899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4138.1875
$ws.Range("I7").Value = 2918.1667
$ws.Range("K7").Value = 2918.1667
$ws.Range("M7").Value = -2806.1667
$ws.Range("H126").Value = 4138.1875
$ws.Range("I126").Value = 2918.1667
$ws.Range("K126").Value = 8754.500100000001
$ws.Range("M126").Value = -6284.500100000001
$ws.Range("H132").Value = 2210.6365
$ws.Range("I132").Value = 1500.6086
$ws.Range("J132").Value = 3843.7
$ws.Range("K132").Value = 4501.825800000001
$ws.Range("L132").Value = 11531.1
$ws.Range("M132").Value = -1971.825800000001
$ws.Range("N132").Value = -16591.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2330.3
$ws.Range("I100").Value = 2329
$ws.Range("J100").Value = 2333.3333
$ws.Range("K100").Value = 4658
$ws.Range("L100").Value = 4666.6666
$ws.Range("M100").Value = -4117
$ws.Range("N100").Value = -5748.6666
$ws.Range("H132").Value = 4506841
$ws.Range("I132").Value = 2742.45
$ws.Range("J132").Value = 9805780
$ws.Range("K132").Value = 8227.349999999999
$ws.Range("L132").Value = 29417340
$ws.Range("M132").Value = -5697.349999999999
$ws.Range("N132").Value = -29422400
